# Update the guest-lecture date shown on the title slide.
# "Text Placeholder 12" (shape id 13) on slide 1 holds the lecture date,
# which moves from the 2023 talk ("20 March 2023") to the 2025 guest
# lecture ("2 April 2025"). Setting .Text on the existing TextRange keeps
# the shape's run-level formatting (e.g. the white "bg1" font colour)
# intact, matching how PowerPoint performs an in-place text edit.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(1)
$dateShape = $slide.Shapes.Item(4)
$dateShape.TextFrame.TextRange.Text = "2 April 2025"
